# ---- Part 1: Update "总计" sheet with new 2022-Q3 row, shifting data down ----
$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

$summary.Cells.Item(2,1).Value = 0
$summary.Cells.Item(2,2).Value = "2022-Q3"
$summary.Cells.Item(2,3).Value = 44
$summary.Cells.Item(2,4).Value = 8.66

$summary.Cells.Item(3,1).Value = 1
$summary.Cells.Item(3,2).Value = "2022-Q2"
$summary.Cells.Item(3,3).Value = 88
$summary.Cells.Item(3,4).Value = 28.32

$summary.Cells.Item(4,1).Value = 2
$summary.Cells.Item(4,2).Value = "2022-Q1"
$summary.Cells.Item(4,3).Value = 63
$summary.Cells.Item(4,4).Value = 32.05

$summary.Cells.Item(5,1).Value = 3
$summary.Cells.Item(5,2).Value = "2021-Q4"
$summary.Cells.Item(5,3).Value = 41
$summary.Cells.Item(5,4).Value = 15.83

$summary.Cells.Item(6,1).Value = 4
$summary.Cells.Item(6,2).Value = "2021-Q3"
$summary.Cells.Item(6,3).Value = 83
$summary.Cells.Item(6,4).Value = 88.45

$summary.Cells.Item(7,1).Value = 5
$summary.Cells.Item(7,2).Value = "2021-Q2"
$summary.Cells.Item(7,3).Value = 48
$summary.Cells.Item(7,4).Value = 23.29

$summary.Cells.Item(8,1).Value = 6
$summary.Cells.Item(8,2).Value = "2021-Q1"
$summary.Cells.Item(8,3).Value = 37
$summary.Cells.Item(8,4).Value = 18.03

$summary.Cells.Item(9,1).Value = 7
$summary.Cells.Item(9,2).Value = "2020-Q4"
$summary.Cells.Item(9,3).Value = 3
$summary.Cells.Item(9,4).Value = 1.44

# Apply matching bold/border/center-top style to the newly created A9 index cell
$rngA9 = $summary.Range("A9")
$rngA9.Font.Bold = $true
$rngA9.HorizontalAlignment = -4108
$rngA9.VerticalAlignment = -4160
$rngA9.Borders.LineStyle = 1

# ---- Part 2: Insert new "2022-Q3" worksheet right after "总计" ----
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q3"
$ws = $newSheet

# Header row
$ws.Cells.Item(1,2).Value = "基金代码"
$ws.Cells.Item(1,3).Value = "基金名称"
$ws.Cells.Item(1,4).Value = "基金规模"
$ws.Cells.Item(1,5).Value = "股票总仓位"
$ws.Cells.Item(1,6).Value = "仓位占比"
$ws.Cells.Item(1,7).Value = "持有市值(亿元)"
$ws.Cells.Item(1,8).Value = "仓位排名"

# Data rows
$ws.Cells.Item(2,1).Value = 0
$ws.Cells.Item(2,2).Value = "'002708"
$ws.Cells.Item(2,3).Value = "大摩健康产业混合A"
$ws.Cells.Item(2,4).Value = "'21.51"
$ws.Cells.Item(2,5).Value = "'92.37"
$ws.Cells.Item(2,6).Value = "'8.82"
$ws.Cells.Item(2,7).Value = "'1.8972"
$ws.Cells.Item(2,8).Value = 4

$ws.Cells.Item(3,1).Value = 1
$ws.Cells.Item(3,2).Value = "'000831"
$ws.Cells.Item(3,3).Value = "工银医疗保健股票"
$ws.Cells.Item(3,4).Value = "'41.69"
$ws.Cells.Item(3,5).Value = "'82.62"
$ws.Cells.Item(3,6).Value = "'4.29"
$ws.Cells.Item(3,7).Value = "'1.7885"
$ws.Cells.Item(3,8).Value = 9

$ws.Cells.Item(4,1).Value = 2
$ws.Cells.Item(4,2).Value = "'001171"
$ws.Cells.Item(4,3).Value = "工银养老产业股票A"
$ws.Cells.Item(4,4).Value = "'25.10"
$ws.Cells.Item(4,5).Value = "'80.17"
$ws.Cells.Item(4,6).Value = "'3.24"
$ws.Cells.Item(4,7).Value = "'0.8132"
$ws.Cells.Item(4,8).Value = 9

$ws.Cells.Item(5,1).Value = 3
$ws.Cells.Item(5,2).Value = "'161616"
$ws.Cells.Item(5,3).Value = "融通医疗保健行业混合A/B"
$ws.Cells.Item(5,4).Value = "'12.82"
$ws.Cells.Item(5,5).Value = "'88.96"
$ws.Cells.Item(5,6).Value = "'3.82"
$ws.Cells.Item(5,7).Value = "'0.4897"
$ws.Cells.Item(5,8).Value = 8

$ws.Cells.Item(6,1).Value = 4
$ws.Cells.Item(6,2).Value = "'014030"
$ws.Cells.Item(6,3).Value = "大摩健康产业混合C"
$ws.Cells.Item(6,4).Value = "'4.27"
$ws.Cells.Item(6,5).Value = "'92.37"
$ws.Cells.Item(6,6).Value = "'8.82"
$ws.Cells.Item(6,7).Value = "'0.3766"
$ws.Cells.Item(6,8).Value = 4

$ws.Cells.Item(7,1).Value = 5
$ws.Cells.Item(7,2).Value = "'010314"
$ws.Cells.Item(7,3).Value = "摩根士丹利华鑫内需增长混合A"
$ws.Cells.Item(7,4).Value = "'4.32"
$ws.Cells.Item(7,5).Value = "'93.24"
$ws.Cells.Item(7,6).Value = "'8.38"
$ws.Cells.Item(7,7).Value = "'0.3620"
$ws.Cells.Item(7,8).Value = 6

$ws.Cells.Item(8,1).Value = 6
$ws.Cells.Item(8,2).Value = "'233007"
$ws.Cells.Item(8,3).Value = "大摩卓越成长混合"
$ws.Cells.Item(8,4).Value = "'4.35"
$ws.Cells.Item(8,5).Value = "'94.21"
$ws.Cells.Item(8,6).Value = "'8.03"
$ws.Cells.Item(8,7).Value = "'0.3493"
$ws.Cells.Item(8,8).Value = 5

$ws.Cells.Item(9,1).Value = 7
$ws.Cells.Item(9,2).Value = "'000945"
$ws.Cells.Item(9,3).Value = "华夏医疗健康混合A"
$ws.Cells.Item(9,4).Value = "'10.20"
$ws.Cells.Item(9,5).Value = "'84.42"
$ws.Cells.Item(9,6).Value = "'2.55"
$ws.Cells.Item(9,7).Value = "'0.2601"
$ws.Cells.Item(9,8).Value = 10

$ws.Cells.Item(10,1).Value = 8
$ws.Cells.Item(10,2).Value = "'012368"
$ws.Cells.Item(10,3).Value = "摩根士丹利华鑫优享臻选六个月持有期混合A"
$ws.Cells.Item(10,4).Value = "'4.42"
$ws.Cells.Item(10,5).Value = "'94.05"
$ws.Cells.Item(10,6).Value = "'5.86"
$ws.Cells.Item(10,7).Value = "'0.2590"
$ws.Cells.Item(10,8).Value = 7

$ws.Cells.Item(11,1).Value = 9
$ws.Cells.Item(11,2).Value = "'000309"
$ws.Cells.Item(11,3).Value = "大摩品质生活精选股票"
$ws.Cells.Item(11,4).Value = "'3.42"
$ws.Cells.Item(11,5).Value = "'94.27"
$ws.Cells.Item(11,6).Value = "'5.88"
$ws.Cells.Item(11,7).Value = "'0.2011"
$ws.Cells.Item(11,8).Value = 9

$ws.Cells.Item(12,1).Value = 10
$ws.Cells.Item(12,2).Value = "'233006"
$ws.Cells.Item(12,3).Value = "大摩领先优势混合"
$ws.Cells.Item(12,4).Value = "'3.47"
$ws.Cells.Item(12,5).Value = "'94.24"
$ws.Cells.Item(12,6).Value = "'5.39"
$ws.Cells.Item(12,7).Value = "'0.1870"
$ws.Cells.Item(12,8).Value = 7

$ws.Cells.Item(13,1).Value = 11
$ws.Cells.Item(13,2).Value = "'004905"
$ws.Cells.Item(13,3).Value = "华泰柏瑞生物医药混合A"
$ws.Cells.Item(13,4).Value = "'4.92"
$ws.Cells.Item(13,5).Value = "'94.10"
$ws.Cells.Item(13,6).Value = "'3.43"
$ws.Cells.Item(13,7).Value = "'0.1688"
$ws.Cells.Item(13,8).Value = 8

$ws.Cells.Item(14,1).Value = 12
$ws.Cells.Item(14,2).Value = "'010322"
$ws.Cells.Item(14,3).Value = "大摩新兴产业股票"
$ws.Cells.Item(14,4).Value = "'1.98"
$ws.Cells.Item(14,5).Value = "'94.19"
$ws.Cells.Item(14,6).Value = "'7.20"
$ws.Cells.Item(14,7).Value = "'0.1426"
$ws.Cells.Item(14,8).Value = 4

$ws.Cells.Item(15,1).Value = 13
$ws.Cells.Item(15,2).Value = "'000523"
$ws.Cells.Item(15,3).Value = "国投瑞银医疗保健混合A"
$ws.Cells.Item(15,4).Value = "'1.95"
$ws.Cells.Item(15,5).Value = "'92.21"
$ws.Cells.Item(15,6).Value = "'7.26"
$ws.Cells.Item(15,7).Value = "'0.1416"
$ws.Cells.Item(15,8).Value = 4

$ws.Cells.Item(16,1).Value = 14
$ws.Cells.Item(16,2).Value = "'013357"
$ws.Cells.Item(16,3).Value = "大摩沪港深精选混合C"
$ws.Cells.Item(16,4).Value = "'1.53"
$ws.Cells.Item(16,5).Value = "'92.27"
$ws.Cells.Item(16,6).Value = "'8.96"
$ws.Cells.Item(16,7).Value = "'0.1371"
$ws.Cells.Item(16,8).Value = 5

$ws.Cells.Item(17,1).Value = 15
$ws.Cells.Item(17,2).Value = "'014867"
$ws.Cells.Item(17,3).Value = "摩根士丹利华鑫优悦安和混合C"
$ws.Cells.Item(17,4).Value = "'1.32"
$ws.Cells.Item(17,5).Value = "'93.41"
$ws.Cells.Item(17,6).Value = "'9.64"
$ws.Cells.Item(17,7).Value = "'0.1272"
$ws.Cells.Item(17,8).Value = 3

$ws.Cells.Item(18,1).Value = 16
$ws.Cells.Item(18,2).Value = "'013349"
$ws.Cells.Item(18,3).Value = "创金合信大健康混合C"
$ws.Cells.Item(18,4).Value = "'1.35"
$ws.Cells.Item(18,5).Value = "'90.70"
$ws.Cells.Item(18,6).Value = "'6.94"
$ws.Cells.Item(18,7).Value = "'0.0937"
$ws.Cells.Item(18,8).Value = 7

$ws.Cells.Item(19,1).Value = 17
$ws.Cells.Item(19,2).Value = "'506009"
$ws.Cells.Item(19,3).Value = "国泰科创板两年定期开放混合"
$ws.Cells.Item(19,4).Value = "'2.05"
$ws.Cells.Item(19,5).Value = "'85.80"
$ws.Cells.Item(19,6).Value = "'4.35"
$ws.Cells.Item(19,7).Value = "'0.0892"
$ws.Cells.Item(19,8).Value = 4

$ws.Cells.Item(20,1).Value = 18
$ws.Cells.Item(20,2).Value = "'163001"
$ws.Cells.Item(20,3).Value = "长信医疗保健行业灵活配置混合（LOF）"
$ws.Cells.Item(20,4).Value = "'1.88"
$ws.Cells.Item(20,5).Value = "'94.04"
$ws.Cells.Item(20,6).Value = "'4.49"
$ws.Cells.Item(20,7).Value = "'0.0844"
$ws.Cells.Item(20,8).Value = 10

$ws.Cells.Item(21,1).Value = 19
$ws.Cells.Item(21,2).Value = "'009893"
$ws.Cells.Item(21,3).Value = "摩根士丹利华鑫优悦安和混合A"
$ws.Cells.Item(21,4).Value = "'0.87"
$ws.Cells.Item(21,5).Value = "'93.41"
$ws.Cells.Item(21,6).Value = "'9.64"
$ws.Cells.Item(21,7).Value = "'0.0839"
$ws.Cells.Item(21,8).Value = 3

$ws.Cells.Item(22,1).Value = 20
$ws.Cells.Item(22,2).Value = "'012238"
$ws.Cells.Item(22,3).Value = "工银养老产业股票C"
$ws.Cells.Item(22,4).Value = "'2.53"
$ws.Cells.Item(22,5).Value = "'80.17"
$ws.Cells.Item(22,6).Value = "'3.24"
$ws.Cells.Item(22,7).Value = "'0.0820"
$ws.Cells.Item(22,8).Value = 9

$ws.Cells.Item(23,1).Value = 21
$ws.Cells.Item(23,2).Value = "'000946"
$ws.Cells.Item(23,3).Value = "华夏医疗健康混合C"
$ws.Cells.Item(23,4).Value = "'2.78"
$ws.Cells.Item(23,5).Value = "'84.42"
$ws.Cells.Item(23,6).Value = "'2.55"
$ws.Cells.Item(23,7).Value = "'0.0709"
$ws.Cells.Item(23,8).Value = 10

$ws.Cells.Item(24,1).Value = 22
$ws.Cells.Item(24,2).Value = "'013356"
$ws.Cells.Item(24,3).Value = "大摩沪港深精选混合A"
$ws.Cells.Item(24,4).Value = "'0.68"
$ws.Cells.Item(24,5).Value = "'92.27"
$ws.Cells.Item(24,6).Value = "'8.96"
$ws.Cells.Item(24,7).Value = "'0.0609"
$ws.Cells.Item(24,8).Value = 5

$ws.Cells.Item(25,1).Value = 23
$ws.Cells.Item(25,2).Value = "'013348"
$ws.Cells.Item(25,3).Value = "创金合信大健康混合A"
$ws.Cells.Item(25,4).Value = "'0.79"
$ws.Cells.Item(25,5).Value = "'90.70"
$ws.Cells.Item(25,6).Value = "'6.94"
$ws.Cells.Item(25,7).Value = "'0.0548"
$ws.Cells.Item(25,8).Value = 7

$ws.Cells.Item(26,1).Value = 24
$ws.Cells.Item(26,2).Value = "'015052"
$ws.Cells.Item(26,3).Value = "东方红医疗升级股票A"
$ws.Cells.Item(26,4).Value = "'1.79"
$ws.Cells.Item(26,5).Value = "'86.95"
$ws.Cells.Item(26,6).Value = "'2.93"
$ws.Cells.Item(26,7).Value = "'0.0524"
$ws.Cells.Item(26,8).Value = 10

$ws.Cells.Item(27,1).Value = 25
$ws.Cells.Item(27,2).Value = "'008923"
$ws.Cells.Item(27,3).Value = "建信医疗健康行业股票A"
$ws.Cells.Item(27,4).Value = "'1.29"
$ws.Cells.Item(27,5).Value = "'86.41"
$ws.Cells.Item(27,6).Value = "'3.12"
$ws.Cells.Item(27,7).Value = "'0.0402"
$ws.Cells.Item(27,8).Value = 9

$ws.Cells.Item(28,1).Value = 26
$ws.Cells.Item(28,2).Value = "'001294"
$ws.Cells.Item(28,3).Value = "新华战略新兴产业灵活配置混合"
$ws.Cells.Item(28,4).Value = "'0.99"
$ws.Cells.Item(28,5).Value = "'93.49"
$ws.Cells.Item(28,6).Value = "'3.68"
$ws.Cells.Item(28,7).Value = "'0.0364"
$ws.Cells.Item(28,8).Value = 6

$ws.Cells.Item(29,1).Value = 27
$ws.Cells.Item(29,2).Value = "'009275"
$ws.Cells.Item(29,3).Value = "融通医疗保健行业混合C"
$ws.Cells.Item(29,4).Value = "'0.92"
$ws.Cells.Item(29,5).Value = "'88.96"
$ws.Cells.Item(29,6).Value = "'3.82"
$ws.Cells.Item(29,7).Value = "'0.0351"
$ws.Cells.Item(29,8).Value = 8

$ws.Cells.Item(30,1).Value = 28
$ws.Cells.Item(30,2).Value = "'014220"
$ws.Cells.Item(30,3).Value = "恒越医疗健康精选混合A"
$ws.Cells.Item(30,4).Value = "'0.72"
$ws.Cells.Item(30,5).Value = "'88.76"
$ws.Cells.Item(30,6).Value = "'4.43"
$ws.Cells.Item(30,7).Value = "'0.0319"
$ws.Cells.Item(30,8).Value = 7

$ws.Cells.Item(31,1).Value = 29
$ws.Cells.Item(31,2).Value = "'005520"
$ws.Cells.Item(31,3).Value = "国投瑞银创新医疗混合"
$ws.Cells.Item(31,4).Value = "'0.41"
$ws.Cells.Item(31,5).Value = "'93.90"
$ws.Cells.Item(31,6).Value = "'7.47"
$ws.Cells.Item(31,7).Value = "'0.0306"
$ws.Cells.Item(31,8).Value = 4

$ws.Cells.Item(32,1).Value = 30
$ws.Cells.Item(32,2).Value = "'008924"
$ws.Cells.Item(32,3).Value = "建信医疗健康行业股票C"
$ws.Cells.Item(32,4).Value = "'0.77"
$ws.Cells.Item(32,5).Value = "'86.41"
$ws.Cells.Item(32,6).Value = "'3.12"
$ws.Cells.Item(32,7).Value = "'0.0240"
$ws.Cells.Item(32,8).Value = 9

$ws.Cells.Item(33,1).Value = 31
$ws.Cells.Item(33,2).Value = "'012496"
$ws.Cells.Item(33,3).Value = "同泰行业优选股票A"
$ws.Cells.Item(33,4).Value = "'0.59"
$ws.Cells.Item(33,5).Value = "'94.78"
$ws.Cells.Item(33,6).Value = "'3.09"
$ws.Cells.Item(33,7).Value = "'0.0182"
$ws.Cells.Item(33,8).Value = 8

$ws.Cells.Item(34,1).Value = 32
$ws.Cells.Item(34,2).Value = "'012369"
$ws.Cells.Item(34,3).Value = "摩根士丹利华鑫优享臻选六个月持有期混合C"
$ws.Cells.Item(34,4).Value = "'0.30"
$ws.Cells.Item(34,5).Value = "'94.05"
$ws.Cells.Item(34,6).Value = "'5.86"
$ws.Cells.Item(34,7).Value = "'0.0176"
$ws.Cells.Item(34,8).Value = 7

$ws.Cells.Item(35,1).Value = 33
$ws.Cells.Item(35,2).Value = "'014221"
$ws.Cells.Item(35,3).Value = "恒越医疗健康精选混合C"
$ws.Cells.Item(35,4).Value = "'0.29"
$ws.Cells.Item(35,5).Value = "'88.76"
$ws.Cells.Item(35,6).Value = "'4.43"
$ws.Cells.Item(35,7).Value = "'0.0128"
$ws.Cells.Item(35,8).Value = 7

$ws.Cells.Item(36,1).Value = 34
$ws.Cells.Item(36,2).Value = "'010031"
$ws.Cells.Item(36,3).Value = "华泰柏瑞生物医药混合C"
$ws.Cells.Item(36,4).Value = "'0.24"
$ws.Cells.Item(36,5).Value = "'94.10"
$ws.Cells.Item(36,6).Value = "'3.43"
$ws.Cells.Item(36,7).Value = "'0.0082"
$ws.Cells.Item(36,8).Value = 8

$ws.Cells.Item(37,1).Value = 35
$ws.Cells.Item(37,2).Value = "'011082"
$ws.Cells.Item(37,3).Value = "国投瑞银医疗保健混合C"
$ws.Cells.Item(37,4).Value = "'0.11"
$ws.Cells.Item(37,5).Value = "'92.21"
$ws.Cells.Item(37,6).Value = "'7.26"
$ws.Cells.Item(37,7).Value = "'0.0080"
$ws.Cells.Item(37,8).Value = 4

$ws.Cells.Item(38,1).Value = 36
$ws.Cells.Item(38,2).Value = "'012497"
$ws.Cells.Item(38,3).Value = "同泰行业优选股票C"
$ws.Cells.Item(38,4).Value = "'0.19"
$ws.Cells.Item(38,5).Value = "'94.78"
$ws.Cells.Item(38,6).Value = "'3.09"
$ws.Cells.Item(38,7).Value = "'0.0059"
$ws.Cells.Item(38,8).Value = 8

$ws.Cells.Item(39,1).Value = 37
$ws.Cells.Item(39,2).Value = "'008842"
$ws.Cells.Item(39,3).Value = "同泰远见灵活配置混合A"
$ws.Cells.Item(39,4).Value = "'0.18"
$ws.Cells.Item(39,5).Value = "'93.90"
$ws.Cells.Item(39,6).Value = "'2.72"
$ws.Cells.Item(39,7).Value = "'0.0049"
$ws.Cells.Item(39,8).Value = 9

$ws.Cells.Item(40,1).Value = 38
$ws.Cells.Item(40,2).Value = "'015053"
$ws.Cells.Item(40,3).Value = "东方红医疗升级股票C"
$ws.Cells.Item(40,4).Value = "'0.13"
$ws.Cells.Item(40,5).Value = "'86.95"
$ws.Cells.Item(40,6).Value = "'2.93"
$ws.Cells.Item(40,7).Value = "'0.0038"
$ws.Cells.Item(40,8).Value = 10

$ws.Cells.Item(41,1).Value = 39
$ws.Cells.Item(41,2).Value = "'014869"
$ws.Cells.Item(41,3).Value = "摩根士丹利华鑫内需增长混合C"
$ws.Cells.Item(41,4).Value = "'0.04"
$ws.Cells.Item(41,5).Value = "'93.24"
$ws.Cells.Item(41,6).Value = "'8.38"
$ws.Cells.Item(41,7).Value = "'0.0034"
$ws.Cells.Item(41,8).Value = 6

$ws.Cells.Item(42,1).Value = 40
$ws.Cells.Item(42,2).Value = "'013154"
$ws.Cells.Item(42,3).Value = "长信医疗保健行业灵活配置混合(LOF)C"
$ws.Cells.Item(42,4).Value = "'0.05"
$ws.Cells.Item(42,5).Value = "'94.04"
$ws.Cells.Item(42,6).Value = "'4.49"
$ws.Cells.Item(42,7).Value = "'0.0022"
$ws.Cells.Item(42,8).Value = 10

$ws.Cells.Item(43,1).Value = 41
$ws.Cells.Item(43,2).Value = "'008843"
$ws.Cells.Item(43,3).Value = "同泰远见灵活配置混合C"
$ws.Cells.Item(43,4).Value = "'0.07"
$ws.Cells.Item(43,5).Value = "'93.90"
$ws.Cells.Item(43,6).Value = "'2.72"
$ws.Cells.Item(43,7).Value = "'0.0019"
$ws.Cells.Item(43,8).Value = 9

$ws.Cells.Item(44,1).Value = 42
$ws.Cells.Item(44,2).Value = "'004724"
$ws.Cells.Item(44,3).Value = "先锋聚元灵活配置混合A"
$ws.Cells.Item(44,4).Value = "'0.04"
$ws.Cells.Item(44,5).Value = "'94.36"
$ws.Cells.Item(44,6).Value = "'2.77"
$ws.Cells.Item(44,7).Value = "'0.0011"
$ws.Cells.Item(44,8).Value = 5

$ws.Cells.Item(45,1).Value = 43
$ws.Cells.Item(45,2).Value = "'004725"
$ws.Cells.Item(45,3).Value = "先锋聚元灵活配置混合C"
$ws.Cells.Item(45,4).Value = "'0.04"
$ws.Cells.Item(45,5).Value = "'94.36"
$ws.Cells.Item(45,6).Value = "'2.77"
$ws.Cells.Item(45,7).Value = "'0.0011"
$ws.Cells.Item(45,8).Value = 5

# Apply header style (bold, centered, top-aligned, bordered) to B1:H1
$rngHeader = $ws.Range("B1:H1")
$rngHeader.Font.Bold = $true
$rngHeader.HorizontalAlignment = -4108
$rngHeader.VerticalAlignment = -4160
$rngHeader.Borders.LineStyle = 1

# Apply index-column style (bold, centered, top-aligned, bordered) to A2:A45
$rngIndex = $ws.Range("A2:A45")
$rngIndex.Font.Bold = $true
$rngIndex.HorizontalAlignment = -4108
$rngIndex.VerticalAlignment = -4160
$rngIndex.Borders.LineStyle = 1

Write-Output "2022-Q3 sheet inserted and 总计 sheet updated successfully."